$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-24 18:54:46"
$wsDeDe.Range("H3").Value = "2016-08-24 18:54:46"
$wsZhCn.Range("H3").Value = "2016-08-24 18:54:41"
$wsZhCn.Range("K3").Value = "2016-08-24 18:54:58"
$wsDeDe.Range("K3").Value = "2016-08-24 18:55:18"
